$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Insert 4 new rows before current row 9 (the "2Y" row), shifting existing
# rows 9-14 down to 13-18.
$ws.Rows("9:12").Insert()

# New futures rows inserted at 9-12
$newData = @(
    @("5M", "SQF26", "FUTURE", 96.3),
    @("0M", "SQQ25", "FUTURE", 95.78749999999999),
    @("2M", "SQV25", "FUTURE", 96.01000000000001),
    @("3M", "SQX25", "FUTURE", 96.125)
)

$r = 9
foreach ($row in $newData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
